$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (C1 text changes "remarks" -> "Expected ")
$ws.Range("C1").Value = "Expected "

# Update B2 email "s@w.com" -> "sw@g.com"
$ws.Range("B2").Value = "sw@g.com"

# Update B5 email "sj@g.com" -> "sk@g.com"
$ws.Range("B5").Value = "sk@g.com"

# Apply bold font + yellow fill to header row A1:C1
# (build the combined style on a scratch cell first, then paste the
# formatting onto the header range so only one new style entry is created)
$scratch = $ws.Range("E1")
$scratch.Font.Bold = $true
$scratch.Interior.Color = 65535
$scratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratch.Clear()

# Update selection to B2
$ws.Range("B2").Select()

# Set page orientation to portrait
$ws.PageSetup.Orientation = 1
